$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 (Tr) value
$ws.Range("B4").Value = 516

# Update row 5 (J3) value
$ws.Range("B5").Value = 411

# Rows 6-16: new labels and values (shifted up from old rows 6-20, some collapsed)
$ws.Range("A6").Value = "Tr-J"
$ws.Range("B6").Value = 338

$ws.Range("A7").Value = "Pg"
$ws.Range("B7").Value = 193

$ws.Range("A8").Value = "K2-Pg"
$ws.Range("B8").Value = 171

$ws.Range("A9").Value = "J2"
$ws.Range("B9").Value = 154

$ws.Range("A10").Value = "J1-J2"
$ws.Range("B10").Value = 109

$ws.Range("A11").Value = "Mz"
$ws.Range("B11").Value = 88

$ws.Range("A12").Value = "J1"
$ws.Range("B12").Value = 69

$ws.Range("A13").Value = "K"
$ws.Range("B13").Value = 62

$ws.Range("A14").Value = "J"
$ws.Range("B14").Value = 46

$ws.Range("A15").Value = "J-K"
$ws.Range("B15").Value = 37

$ws.Range("A16").Value = "J2-J3"
$ws.Range("B16").Value = 35

# Remove old rows 17-20 (now empty tail of the table)
$ws.Range("A17:B20").EntireRow.Delete()
